$d = $word.ActiveDocument

$d.Content.Find.Execute("431×8=3448", $true, $false, $false, $false, $false, $true, 1, $false, "549×4=2196", 2) | Out-Null
$d.Content.Find.Execute("905×8=7240", $true, $false, $false, $false, $false, $true, 1, $false, "852×7=5964", 2) | Out-Null
$d.Content.Find.Execute("528×2=1056", $true, $false, $false, $false, $false, $true, 1, $false, "793×5=3965", 2) | Out-Null
$d.Content.Find.Execute("795×8=6360", $true, $false, $false, $false, $false, $true, 1, $false, "883×2=1766", 2) | Out-Null
$d.Content.Find.Execute("384×8=3072", $true, $false, $false, $false, $false, $true, 1, $false, "231×6=1386", 2) | Out-Null
$d.Content.Find.Execute("847×2=1694", $true, $false, $false, $false, $false, $true, 1, $false, "185×3=555", 2) | Out-Null
$d.Content.Find.Execute("718×2=1436", $true, $false, $false, $false, $false, $true, 1, $false, "594×7=4158", 2) | Out-Null
$d.Content.Find.Execute("599×9=5391", $true, $false, $false, $false, $false, $true, 1, $false, "191×8=1528", 2) | Out-Null
$d.Content.Find.Execute("515×8=4120", $true, $false, $false, $false, $false, $true, 1, $false, "809×8=6472", 2) | Out-Null
$d.Content.Find.Execute("477×2=954", $true, $false, $false, $false, $false, $true, 1, $false, "426×4=1704", 2) | Out-Null
$d.Content.Find.Execute("670×4=2680", $true, $false, $false, $false, $false, $true, 1, $false, "778×3=2334", 2) | Out-Null
$d.Content.Find.Execute("970×5=4850", $true, $false, $false, $false, $false, $true, 1, $false, "381×9=3429", 2) | Out-Null
$d.Content.Find.Execute("868×5=4340", $true, $false, $false, $false, $false, $true, 1, $false, "875×4=3500", 2) | Out-Null
$d.Content.Find.Execute("738×2=1476", $true, $false, $false, $false, $false, $true, 1, $false, "174×7=1218", 2) | Out-Null
$d.Content.Find.Execute("357×6=2142", $true, $false, $false, $false, $false, $true, 1, $false, "113×6=678", 2) | Out-Null
$d.Content.Find.Execute("922×9=8298", $true, $false, $false, $false, $false, $true, 1, $false, "576×7=4032", 2) | Out-Null
$d.Content.Find.Execute("650×6=3900", $true, $false, $false, $false, $false, $true, 1, $false, "656×3=1968", 2) | Out-Null
$d.Content.Find.Execute("362×6=2172", $true, $false, $false, $false, $false, $true, 1, $false, "849×7=5943", 2) | Out-Null
$d.Content.Find.Execute("757×8=6056", $true, $false, $false, $false, $false, $true, 1, $false, "679×6=4074", 2) | Out-Null
$d.Content.Find.Execute("475×8=3800", $true, $false, $false, $false, $false, $true, 1, $false, "261×8=2088", 2) | Out-Null
$d.Content.Find.Execute("818×3=2454", $true, $false, $false, $false, $false, $true, 1, $false, "543×2=1086", 2) | Out-Null
$d.Content.Find.Execute("611×3=1833", $true, $false, $false, $false, $false, $true, 1, $false, "876×5=4380", 2) | Out-Null
$d.Content.Find.Execute("559×2=1118", $true, $false, $false, $false, $false, $true, 1, $false, "640×6=3840", 2) | Out-Null
$d.Content.Find.Execute("800×9=7200", $true, $false, $false, $false, $false, $true, 1, $false, "730×3=2190", 2) | Out-Null
$d.Content.Find.Execute("485×8=3880", $true, $false, $false, $false, $false, $true, 1, $false, "329×5=1645", 2) | Out-Null
